$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Add the two new columns (Mid Paper 1 / Mid Paper 2) on both sheets ---
foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("F1").Value = "Mid Paper 1"
    $ws.Range("G1").Value = "Mid Paper 2"

    # Mid-term marks for the first student mirror the final paper marks
    $ws.Range("F2").Value = 67
    $ws.Range("G2").Value = 78
}

# --- Update the active sheet / selection so "Senior Five" becomes the selected tab ---
[void]$ws1.Range("F1").Select()
[void]$ws2.Activate()
[void]$ws2.Range("F10").Select()
